$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 64:77 (no longer present in the updated dataset)
$ws.Range("A64:A77").EntireRow.Delete()

# Update case-count values for rows 2 through 63
$values = @(
    3599, 773, 701, 563, 338, 212, 153, 89, 75, 73, 68, 33, 30, 29, 22, 22, 21, 21, 19, 16, 15, 12, 11, 10, 9, 8, 8, 8, 7, 6, 6, 6, 5, 5, 4, 3, 3, 3, 3, 3, 3, 3, 3, 2, 2, 2, 2, 2, 2, 2, 2, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
